# ============================================================================
# Adds two new worksheets (qrCode, passwordEncryption) with the "dak register"
# data used to build per-row QR payload strings, plus a small password table
# keyed by PDF file name. Mirrors the authoring commit that introduced
# password-protected PDF export driven from this workbook.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Set sheet1's selection while it is still the active sheet (Range.Select()
# on another worksheet would otherwise force-activate it first).
$ws1.Range("F1:F2").Select()

# ----------------------------------------------------------------------
# 1. New sheet "qrCode"
# ----------------------------------------------------------------------
$wsQr = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsQr.Name = "qrCode"

# Header row
$wsQr.Range("A1").Value = "Sl_No"
$wsQr.Range("B1").Value = "Control"
$wsQr.Range("C1").Value = "Dak No"
$wsQr.Range("D1").Value = "Dak Date"
$wsQr.Range("E1").Value = "Salutation"
$wsQr.Range("F1").Value = "Name"
$wsQr.Range("G1").Value = "Employee ID"
$wsQr.Range("H1").Value = "From"
$wsQr.Range("I1").Value = "To"
$wsQr.Range("J1").Value = "Conduct"
$wsQr.Range("K1").Value = "Department"
$wsQr.Range("L1").Value = "Designation"
$wsQr.Range("M1").Value = "QR Code"

# Row 2
$wsQr.Range("A2").Value = 1
$wsQr.Range("B2").Value = 1
$wsQr.Range("C2").Value = 34543
$wsQr.Range("D2").NumberFormat = "@"
$wsQr.Range("D2").Value = "10-06-2021"
$wsQr.Range("E2").Value = "Mr."
$wsQr.Range("F2").Value = "Umakanta Pattanaik"
$wsQr.Range("G2").Value = 10200
$wsQr.Range("H2").NumberFormat = "@"
$wsQr.Range("H2").Value = "10-06-2021"
$wsQr.Range("I2").NumberFormat = "@"
$wsQr.Range("I2").Value = "10-06-2021"
$wsQr.Range("J2").Value = "Good"
$wsQr.Range("K2").Value = "Development"
$wsQr.Range("L2").Value = "Field credit Officer"
$wsQr.Range("M2").Formula = '=$F$1&"-"&F2&","&$G$1&"-"&G2&","&$H$1&"-"&H2&","&$I$1&"-"&I2&","&$J$1&"-"&J2'

# Row 3
$wsQr.Range("A3").Value = 2
$wsQr.Range("B3").Value = 1
$wsQr.Range("C3").Value = 23454
$wsQr.Range("D3").NumberFormat = "@"
$wsQr.Range("D3").Value = "10-06-2021"
$wsQr.Range("E3").Value = "Mr."
$wsQr.Range("F3").Value = "Purna Chandra Pattanaik"
$wsQr.Range("G3").Value = 11298
$wsQr.Range("H3").NumberFormat = "@"
$wsQr.Range("H3").Value = "10-06-2021"
$wsQr.Range("I3").NumberFormat = "@"
$wsQr.Range("I3").Value = "10-06-2021"
$wsQr.Range("J3").Value = "Good"
$wsQr.Range("K3").Value = "Sales"
$wsQr.Range("L3").Value = "Field credit Officer"
$wsQr.Range("M3").Formula = '=$F$1&"-"&F3&","&$G$1&"-"&G3&","&$H$1&"-"&H3&","&$I$1&"-"&I3&","&$J$1&"-"&J3'

# Formatting: thin border + centred alignment on the used range
$wsQr.Range("A1:M3").Borders.LineStyle = 1
$wsQr.Range("A1:M3").HorizontalAlignment = -4108
$wsQr.Range("A1:M3").VerticalAlignment = -4108
$wsQr.Range("A1:M3").WrapText = $true
$wsQr.Range("A2:M3").WrapText = $false
$wsQr.Range("D2:D3").WrapText = $false
$wsQr.Range("H2:I3").WrapText = $false

# Placeholder empty cells (H4:I11) kept from the original authoring, left blank
for ($r = 4; $r -le 11; $r++) {
    $wsQr.Cells.Item($r, 8).Value = $null
    $wsQr.Cells.Item($r, 9).Value = $null
}

# Data validation: Salutation list
$wsQr.Range("E2:E3").Validation.Add(3, 1, 1, '"Mr.,Ms."')
$wsQr.Range("E2:E3").Validation.InCellDropdown = $true
$wsQr.Range("E2:E3").Validation.IgnoreBlank = $true

# Column widths (approximate best-fit)
$wsQr.Columns.Item(1).ColumnWidth = 4.7
$wsQr.Columns.Item(2).ColumnWidth = 6.17
$wsQr.Columns.Item(3).ColumnWidth = 5.98
$wsQr.Columns.Item(4).ColumnWidth = 9.26
$wsQr.Columns.Item(5).ColumnWidth = 8.35
$wsQr.Columns.Item(6).ColumnWidth = 20.89
$wsQr.Columns.Item(7).ColumnWidth = 10.44
$wsQr.Columns.Item(8).ColumnWidth = 9.26
$wsQr.Columns.Item(9).ColumnWidth = 9.26
$wsQr.Columns.Item(10).ColumnWidth = 6.89
$wsQr.Columns.Item(11).ColumnWidth = 11.17
$wsQr.Columns.Item(12).ColumnWidth = 15.26
$wsQr.Columns.Item(13).ColumnWidth = 35.35

$wsQr.Range("A1:C3").Select()

# ----------------------------------------------------------------------
# 2. New sheet "passwordEncryption"
# ----------------------------------------------------------------------
$wsPw = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsPw.Name = "passwordEncryption"

$wsPw.Range("A1").Value = "Sl_No"
$wsPw.Range("B1").Value = "Control"
$wsPw.Range("C1").Value = "PDFName"
$wsPw.Range("D1").Value = "Password"

$wsPw.Range("A2").Value = 1
$wsPw.Range("B2").Value = 1
$wsPw.Range("C2").NumberFormat = "@"
$wsPw.Range("C2").Value = "'2345"
$wsPw.Range("D2").Value = "password"

$wsPw.Range("A3").Value = 2
$wsPw.Range("B3").Value = 0
$wsPw.Range("C3").Value = "'7209"
$wsPw.Range("D3").Value = "password"

$wsPw.Range("A4").Value = 3
$wsPw.Range("B4").Value = 1
$wsPw.Range("C4").Value = "'22345"
$wsPw.Range("D4").Value = "password"

$wsPw.Range("A1:B4").Borders.LineStyle = 1
$wsPw.Range("A1:B4").HorizontalAlignment = -4108
$wsPw.Range("A1:B4").VerticalAlignment = -4108

$wsPw.Range("C1:D4").Borders.LineStyle = 1

$wsPw.Columns.Item(3).ColumnWidth = 19.71
$wsPw.Columns.Item(4).ColumnWidth = 9.53

$wsPw.Range("B4").Select()
$wsPw.Activate()
